{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document is a daily \"two-digit \u00f7 one-digit\" division worksheet: a date\n// heading paragraph followed by a table of arithmetic problems laid out as\n// \"A\u00f7B=C, D\" (quotient C, remainder D). This edit swaps the date and every\n// problem/answer string for a new day's set. Formatting (fonts, sizes,\n// alignment, table structure) is left untouched -- only the text runs'\n// characters change, via exact literal find/replace.\n//\n// Every \"old\" string below is unique in the document and none of the new\n// strings collide with any other pair's old string, so a plain sequential\n// search+replace is safe and order-independent.\nconst replacements = [\n  [\"2024-01-17 Wednesday\", \"2024-01-18 Thursday\"],\n  [\"89\u00f75=17, 4\", \"70\u00f72=35, 0\"],\n  [\"33\u00f74=8, 1\", \"75\u00f79=8, 3\"],\n  [\"96\u00f78=12, 0\", \"20\u00f79=2, 2\"],\n  [\"74\u00f73=24, 2\", \"10\u00f76=1, 4\"],\n  [\"50\u00f78=6, 2\", \"41\u00f74=10, 1\"],\n  [\"89\u00f74=22, 1\", \"37\u00f74=9, 1\"],\n  [\"44\u00f74=11, 0\", \"79\u00f76=13, 1\"],\n  [\"22\u00f73=7, 1\", \"65\u00f78=8, 1\"],\n  [\"31\u00f76=5, 1\", \"46\u00f77=6, 4\"],\n  [\"60\u00f79=6, 6\", \"49\u00f77=7, 0\"],\n  [\"10\u00f73=3, 1\", \"52\u00f79=5, 7\"],\n  [\"53\u00f78=6, 5\", \"28\u00f77=4, 0\"],\n  [\"24\u00f79=2, 6\", \"86\u00f78=10, 6\"],\n  [\"17\u00f73=5, 2\", \"51\u00f77=7, 2\"],\n  [\"92\u00f74=23, 0\", \"93\u00f77=13, 2\"],\n  [\"73\u00f76=12, 1\", \"50\u00f75=10, 0\"],\n  [\"94\u00f76=15, 4\", \"23\u00f75=4, 3\"],\n  [\"88\u00f78=11, 0\", \"48\u00f78=6, 0\"],\n  [\"32\u00f73=10, 2\", \"34\u00f76=5, 4\"],\n  [\"92\u00f78=11, 4\", \"75\u00f76=12, 3\"],\n  [\"93\u00f73=31, 0\", \"81\u00f78=10, 1\"],\n  [\"41\u00f79=4, 5\", \"95\u00f72=47, 1\"],\n  [\"30\u00f72=15, 0\", \"39\u00f77=5, 4\"],\n  [\"56\u00f76=9, 2\", \"24\u00f76=4, 0\"],\n  [\"99\u00f74=24, 3\", \"95\u00f73=31, 2\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // matchCase keeps \"\u00f7\"/digits literal; matchWholeWord is irrelevant here\n  // since we always match the full \"A\u00f7B=C, D\" string.\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace every match (there is exactly one per pair in this document).\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# PowerShell / Word COM interop edit script.\n# The document is open as $word.ActiveDocument ($d below).\n#\n# The document is a daily \"two-digit / one-digit division\" worksheet: a date\n# heading paragraph followed by a table of arithmetic problems\n# (\"A/B=C, D\"). This script swaps the date and every problem/answer string\n# for a new day's set via Find/Replace over the whole document body,\n# leaving all formatting untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-01-17 Wednesday', '2024-01-18 Thursday'),\n    @('89\u00f75=17, 4', '70\u00f72=35, 0'),\n    @('33\u00f74=8, 1', '75\u00f79=8, 3'),\n    @('96\u00f78=12, 0', '20\u00f79=2, 2'),\n    @('74\u00f73=24, 2', '10\u00f76=1, 4'),\n    @('50\u00f78=6, 2', '41\u00f74=10, 1'),\n    @('89\u00f74=22, 1', '37\u00f74=9, 1'),\n    @('44\u00f74=11, 0', '79\u00f76=13, 1'),\n    @('22\u00f73=7, 1', '65\u00f78=8, 1'),\n    @('31\u00f76=5, 1', '46\u00f77=6, 4'),\n    @('60\u00f79=6, 6', '49\u00f77=7, 0'),\n    @('10\u00f73=3, 1', '52\u00f79=5, 7'),\n    @('53\u00f78=6, 5', '28\u00f77=4, 0'),\n    @('24\u00f79=2, 6', '86\u00f78=10, 6'),\n    @('17\u00f73=5, 2', '51\u00f77=7, 2'),\n    @('92\u00f74=23, 0', '93\u00f77=13, 2'),\n    @('73\u00f76=12, 1', '50\u00f75=10, 0'),\n    @('94\u00f76=15, 4', '23\u00f75=4, 3'),\n    @('88\u00f78=11, 0', '48\u00f78=6, 0'),\n    @('32\u00f73=10, 2', '34\u00f76=5, 4'),\n    @('92\u00f78=11, 4', '75\u00f76=12, 3'),\n    @('93\u00f73=31, 0', '81\u00f78=10, 1'),\n    @('41\u00f79=4, 5', '95\u00f72=47, 1'),\n    @('30\u00f72=15, 0', '39\u00f77=5, 4'),\n    @('56\u00f76=9, 2', '24\u00f76=4, 0'),\n    @('99\u00f74=24, 3', '95\u00f73=31, 2'),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n"}
